# "watch list script fixes"
#
# Sheet "Test Cases" (first sheet): the "Runmode" column (D) is flipped
# from "Y" to a brand new shared string "N" for (almost) every data row,
# and the "Results" column (E) has a few SKIP/PASS swaps. Rows 12 and 55
# are deliberately left untouched (still "Y"), matching the source diff.
# D70/D71 additionally need their direct formatting brought in line with
# the rest of the column (they were using the plain border-only style
# instead of the wrapped/top-aligned one every other D cell uses).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 2-71 excluding 12, 55 and 58 -- those three keep "Y" in the
# Runmode column untouched. Every other row (including 69) flips to "N".
$runmodeRows = 2..71 | Where-Object { $_ -notin @(12, 55, 58) }

# D70/D71 currently carry a different direct format (s="1") than the rest
# of the column (s="18"). Bring their formatting into line by copying it
# from an already-correct cell before writing the new value.
$ws.Range("D2").Copy()
$ws.Range("D70").PasteSpecial(-4122)
$ws.Range("D71").PasteSpecial(-4122)
$excel.CutCopyMode = $false

foreach ($r in $runmodeRows) {
    $ws.Range("D$r").Value = "N"
}

# Results column fixes
$ws.Range("E58").Value = "PASS"
$ws.Range("E69").Value = "SKIP"
$ws.Range("E70").Value = "SKIP"
$ws.Range("E71").Value = "SKIP"

# Move the saved selection/cursor to D12 (previously B3), and drop the
# scrolled-away top-left anchor so the view opens back at the top.
$ws.Range("D12").Select()
